$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 254.18182
$ws.Range("J2").Value = 534
$ws.Range("L2").Value = 534
$ws.Range("N2").Value = -760
$ws.Range("H28").Value = 606.7143
$ws.Range("I28").Value = 407.4737
$ws.Range("J28").Value = 2499.5
$ws.Range("K28").Value = 407.4737
$ws.Range("L28").Value = 2499.5
$ws.Range("M28").Value = 77.52629999999999
$ws.Range("N28").Value = -3469.5
$ws.Range("H94").Value = 2930.889
$ws.Range("I94").Value = 2930.889
$ws.Range("K94").Value = 2930.889
$ws.Range("M94").Value = -2479.889
$ws.Range("H129").Value = 2013.9445
$ws.Range("I129").Value = 1382.2858
$ws.Range("K129").Value = 4146.857400000001
$ws.Range("M129").Value = 853.1425999999992
$ws.Range("H141").Value = 2801
$ws.Range("J141").Value = 5205
$ws.Range("L141").Value = 15615
$ws.Range("N141").Value = -25975

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4005.2263
$ws.Range("I32").Value = 4184.9775
$ws.Range("K32").Value = 4184.9775
$ws.Range("M32").Value = -3897.9775
$ws.Range("H61").Value = 3605.4473
$ws.Range("I61").Value = 2293.25
$ws.Range("J61").Value = 7279.6
$ws.Range("K61").Value = 2293.25
$ws.Range("L61").Value = 7279.6
$ws.Range("M61").Value = -2081.25
$ws.Range("N61").Value = -7703.6
$ws.Range("H88").Value = 1322.6666
$ws.Range("J88").Value = 1343.3
$ws.Range("L88").Value = 1343.3
$ws.Range("N88").Value = -2155.3
$ws.Range("H91").Value = 1322.6666
$ws.Range("J91").Value = 1343.3
$ws.Range("L91").Value = 1343.3
$ws.Range("N91").Value = -4151.3
$ws.Range("H92").Value = 67958.336
$ws.Range("J92").Value = 67958.336
$ws.Range("L92").Value = 67958.336
$ws.Range("N92").Value = -72950.336
$ws.Range("H136").Value = 3605.4473
$ws.Range("I136").Value = 2293.25
$ws.Range("J136").Value = 7279.6
$ws.Range("K136").Value = 6879.75
$ws.Range("L136").Value = 21838.8
$ws.Range("M136").Value = -4329.75
$ws.Range("N136").Value = -26938.8
$ws.Range("H139").Value = 83571.21000000001
$ws.Range("J139").Value = 83571.21000000001
$ws.Range("L139").Value = 83571.21000000001
$ws.Range("N139").Value = -93851.21000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4752.921
$ws.Range("I134").Value = 2444.28
$ws.Range("K134").Value = 7332.84
$ws.Range("M134").Value = -4797.84

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 586.7692
$ws.Range("I16").Value = 586.7692
$ws.Range("K16").Value = 586.7692
$ws.Range("M16").Value = -299.7692
$ws.Range("H22").Value = 1522.5
$ws.Range("I22").Value = 801.5833
$ws.Range("J22").Value = 2964.3333
$ws.Range("K22").Value = 801.5833
$ws.Range("L22").Value = 2964.3333
$ws.Range("M22").Value = -451.5833
$ws.Range("N22").Value = -3664.3333
$ws.Range("H28").Value = 5643
$ws.Range("J28").Value = 5643
$ws.Range("L28").Value = 5643
$ws.Range("N28").Value = -6133
$ws.Range("H31").Value = 4487.1
$ws.Range("I31").Value = 2979.1904
$ws.Range("K31").Value = 2979.1904
$ws.Range("M31").Value = -2684.1904
$ws.Range("H34").Value = 4487.1
$ws.Range("I34").Value = 2979.1904
$ws.Range("K34").Value = 2979.1904
$ws.Range("M34").Value = -2777.1904
$ws.Range("H92").Value = 70980.60000000001
$ws.Range("J92").Value = 70980.60000000001
$ws.Range("L92").Value = 70980.60000000001
$ws.Range("N92").Value = -75972.60000000001
$ws.Range("H113").Value = 586.7692
$ws.Range("I113").Value = 586.7692
$ws.Range("K113").Value = 586.7692
$ws.Range("M113").Value = 1583.2308
$ws.Range("H122").Value = 3003.6428
$ws.Range("I122").Value = 1334.7142
$ws.Range("J122").Value = 4672.5713
$ws.Range("K122").Value = 4004.1426
$ws.Range("L122").Value = 14017.7139
$ws.Range("M122").Value = -1554.1426
$ws.Range("N122").Value = -18917.7139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 10
$ws.Range("I6").Value = 10
$ws.Range("K6").Value = 30
$ws.Range("M6").Value = 83
$ws.Range("H7").Value = 3153991.2
$ws.Range("I7").Value = 2000184.8
$ws.Range("J7").Value = 7000013.5
$ws.Range("K7").Value = 6000554.4
$ws.Range("L7").Value = 21000040.5
$ws.Range("M7").Value = -6000442.4
$ws.Range("N7").Value = -21000264.5
$ws.Range("H33").Value = 171.33333
$ws.Range("I33").Value = 169
$ws.Range("K33").Value = 1014
$ws.Range("M33").Value = -731
$ws.Range("H120").Value = 27024
$ws.Range("I120").Value = 14015
$ws.Range("K120").Value = 42045
$ws.Range("M120").Value = -37207
$ws.Range("H129").Value = 2009.8334
$ws.Range("I129").Value = 1010.625
$ws.Range("J129").Value = 4008.25
$ws.Range("K129").Value = 3031.875
$ws.Range("L129").Value = 12024.75
$ws.Range("M129").Value = 1968.125
$ws.Range("N129").Value = -22024.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8773.429
$ws.Range("I70").Value = 8581.200000000001
$ws.Range("K70").Value = 8581.200000000001
$ws.Range("M70").Value = -8311.200000000001
$ws.Range("H73").Value = 8773.429
$ws.Range("I73").Value = 8581.200000000001
$ws.Range("K73").Value = 8581.200000000001
$ws.Range("M73").Value = -7645.200000000001
$ws.Range("H98").Value = 74006
$ws.Range("J98").Value = 74006
$ws.Range("L98").Value = 74006
$ws.Range("N98").Value = -79996
$ws.Range("H102").Value = 1897.5333
$ws.Range("I102").Value = 1121.9584
$ws.Range("K102").Value = 1121.9584
$ws.Range("M102").Value = 500.0416
$ws.Range("H123").Value = 75334
$ws.Range("J123").Value = 75334
$ws.Range("L123").Value = 75334
$ws.Range("N123").Value = -80234
$ws.Range("H132").Value = 1322.3529
$ws.Range("I132").Value = 1150.4615
$ws.Range("K132").Value = 3451.3845
$ws.Range("M132").Value = -921.3844999999997

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4666.5557
$ws.Range("I68").Value = 2999.75
$ws.Range("K68").Value = 2999.75
$ws.Range("M68").Value = -2250.75
$ws.Range("H71").Value = 4666.5557
$ws.Range("I71").Value = 2999.75
$ws.Range("K71").Value = 14998.75
$ws.Range("M71").Value = -11254.75
$ws.Range("H93").Value = 3026.111
$ws.Range("I93").Value = 3184.75
$ws.Range("J93").Value = 2899.2
$ws.Range("K93").Value = 3184.75
$ws.Range("L93").Value = 2899.2
$ws.Range("M93").Value = -1936.75
$ws.Range("N93").Value = -5395.2
$ws.Range("H94").Value = 67250
$ws.Range("J94").Value = 67250
$ws.Range("L94").Value = 67250
$ws.Range("N94").Value = -68602
$ws.Range("H100").Value = 6362.207
$ws.Range("I100").Value = 2150.2222
$ws.Range("J100").Value = 13254.546
$ws.Range("K100").Value = 2150.2222
$ws.Range("L100").Value = 13254.546
$ws.Range("M100").Value = -1609.2222
$ws.Range("N100").Value = -14336.546
$ws.Range("H122").Value = 3329.5435
$ws.Range("I122").Value = 3101.8948
$ws.Range("K122").Value = 9305.6844
$ws.Range("M122").Value = -6855.6844

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1110
$ws.Range("J81").Value = 1041.5
$ws.Range("L81").Value = 2083
$ws.Range("N81").Value = -4205
$ws.Range("H84").Value = 1110
$ws.Range("J84").Value = 1041.5
$ws.Range("L84").Value = 10415
$ws.Range("N84").Value = -21023
$ws.Range("H113").Value = 2206.0967
$ws.Range("J113").Value = 3496.2856
$ws.Range("L113").Value = 10488.8568
$ws.Range("N113").Value = -14828.8568
$ws.Range("H126").Value = 4085.3333
$ws.Range("I126").Value = 4109.2666
$ws.Range("J126").Value = 3965.6667
$ws.Range("K126").Value = 12327.7998
$ws.Range("L126").Value = 11897.0001
$ws.Range("M126").Value = -9857.799800000001
$ws.Range("N126").Value = -16837.0001
$ws.Range("H132").Value = 5010.595
$ws.Range("I132").Value = 4185.8057
$ws.Range("K132").Value = 12557.4171
$ws.Range("M132").Value = -10027.4171
